$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.528.35"
$ws.Range("E2").Value = "  -0.40%  "

$ws.Range("D3").Value = "3.786.49"
$ws.Range("E3").Value = "  +0.68%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "614.38"
$ws.Range("E5").Value = "  -1.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.19"
$ws.Range("E6").Value = "  -1.97%  "

$ws.Range("D7").Value = "3.781.52"
$ws.Range("E7").Value = "  +0.69%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("E9").Value = "  -0.95%  "

$ws.Range("E10").Value = "  -1.87%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.43"
$ws.Range("E11").Value = "  +1.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.484"
$ws.Range("E12").Value = "  -1.39%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.85"
$ws.Range("E13").Value = "  -3.68%  "

$ws.Range("E14").Value = "  -2.29%  "

$ws.Range("D15").Value = "4.416.87"
$ws.Range("E15").Value = "  +0.64%  "

$ws.Range("D16").Value = "3.784.72"
$ws.Range("E16").Value = "  +0.69%  "

$ws.Range("D17").Value = "69.599.00"
$ws.Range("E17").Value = "  -0.42%  "

$ws.Range("E18").Value = "  -1.07%  "

$ws.Range("E19").Value = "  -3.67%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "509.22"
$ws.Range("E20").Value = "  +0.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.57"
$ws.Range("E21").Value = "  -1.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.60"
$ws.Range("E22").Value = "  +0.06%  "

$ws.Range("E23").Value = "  +0.71%  "

$ws.Range("E24").Value = "  -1.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.33"
$ws.Range("E25").Value = "  -1.15%  "

$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000142"
$ws.Range("E26").Value = "  +3.74%  "

$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.87"
$ws.Range("E27").Value = "  -2.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.56"
$ws.Range("E28").Value = "  -5.23%  "

$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("E30").Value = "  +3.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.51"
$ws.Range("E31").Value = "  -0.18%  "

$ws.Range("E32").Value = "  +2.74%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.26"
$ws.Range("E33").Value = "  +0.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.114"
$ws.Range("E34").Value = "  -1.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("E36").Value = "  -1.62%  "

$ws.Range("E37").Value = "  -1.66%  "

$ws.Range("E38").Value = "  +6.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "482.06"
$ws.Range("E39").Value = "  +12.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.339"
$ws.Range("E40").Value = "  +0.42%  "

$ws.Range("E41").Value = "  -2.67%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "49.76"

$ws.Range("E43").Value = "  +4.60%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.14"
$ws.Range("E44").Value = "  -3.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.56"
$ws.Range("E45").Value = "  -2.07%  "

$ws.Range("D46").Value = "2.942.66"
$ws.Range("E46").Value = "  -2.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0363"
$ws.Range("E47").Value = "  -0.82%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.39"
$ws.Range("E48").Value = "  -0.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "139.40"
$ws.Range("E49").Value = "  +1.80%  "

$ws.Range("E50").Value = "  +0.05%  "

$ws.Range("E51").Value = "  -1.62%  "
